# Updates cryptos list: prices (col D) and volume% (col E) for rows 2-51,
# plus a coin swap (Aave -> Tezos) on row 51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.768.10"
$ws.Range("E2").Value = "  +5.66%  "

$ws.Range("D3").Value = "1.705.43"
$ws.Range("E3").Value = "  +3.54%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3685"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.19"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3312"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07349"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.194"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.867"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.55%  "

$ws.Range("D16").Value = "1.698.84"
$ws.Range("E16").Value = "  +3.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001068"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06611"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9989"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.057"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.08%  "

$ws.Range("D24").Value = "25.750.02"
$ws.Range("E24").Value = "  +5.66%  "

$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("E26").Value = "  +7.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.307"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.35%  "

$ws.Range("D30").Value = "1.888.75"
$ws.Range("E30").Value = "  +3.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.098"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.950"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08505"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.692"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.338"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.273"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06204"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.535"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2120"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02249"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +17.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6114"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9994"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.846"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5834"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07215"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.72%  "

$ws.Range("B51").Value = "Tezos"
$ws.Range("C51").Value = "https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.82%  "
